$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the existing row 194 (so the old
# rows 194..204 shift down to 196..206), mirroring a new weekly price
# update for "Sandia" (Primera / Segunda) getting prepended.
$ws.Rows.Item(194).Insert()
$ws.Rows.Item(194).Insert()

# New row 194: Terminal Hortofrutícola Agro Chillán - Sandia - Primera
$ws.Cells.Item(194, 1).Value = 7
$ws.Cells.Item(194, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(194, 3).Value = "Ñuble"
$ws.Cells.Item(194, 4).Value = 44918
$ws.Cells.Item(194, 5).Value = 16
$ws.Cells.Item(194, 6).Value = 100112028
$ws.Cells.Item(194, 7).Value = "Sandia"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 1200
$ws.Cells.Item(194, 11).Value = 2800
$ws.Cells.Item(194, 12).Value = 3000
$ws.Cells.Item(194, 13).Value = 2900
$ws.Cells.Item(194, 14).Value = "$/unidad"
$ws.Cells.Item(194, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(194, 16).Value = 2900
$ws.Cells.Item(194, 17).Value = 1
$ws.Cells.Item(194, 18).Value = "Hortaliza"

# New row 195: Terminal Hortofrutícola Agro Chillán - Sandia - Segunda
$ws.Cells.Item(195, 1).Value = 7
$ws.Cells.Item(195, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(195, 3).Value = "Ñuble"
$ws.Cells.Item(195, 4).Value = 44918
$ws.Cells.Item(195, 5).Value = 16
$ws.Cells.Item(195, 6).Value = 100112028
$ws.Cells.Item(195, 7).Value = "Sandia"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Segunda"
$ws.Cells.Item(195, 10).Value = 500
$ws.Cells.Item(195, 11).Value = 2500
$ws.Cells.Item(195, 12).Value = 2500
$ws.Cells.Item(195, 13).Value = 2500
$ws.Cells.Item(195, 14).Value = "$/unidad"
$ws.Cells.Item(195, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(195, 16).Value = 2500
$ws.Cells.Item(195, 17).Value = 1
$ws.Cells.Item(195, 18).Value = "Hortaliza"
